# Change "Device" -> "Apparatus" throughout the workbook.
#
# The commit renames the "Device" worksheet to "Apparatus" and updates the
# three text cells on that sheet that referred to "device(s)". It also
# leaves that sheet as the active/selected tab (with A2 selected), whereas
# previously the "Bus" sheet was the selected tab.

$wb = $excel.ActiveWorkbook

# Rename the "Device" sheet to "Apparatus".
$wsApparatus = $wb.Worksheets.Item("Device")
$wsApparatus.Name = "Apparatus"

# Update the sheet-description / header text from "device" to "apparatus".
$wsApparatus.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."
$wsApparatus.Range("B2").Value = "Apparatus type"
$wsApparatus.Range("C2").Value = "Apparatus parameters"

# Make "Apparatus" the active sheet/tab with cell A2 selected (previously
# "Bus" was the active tab and Device's selection was G6).
[void]$wsApparatus.Activate()
[void]$wsApparatus.Range("A2").Select()
